$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 74; existing rows 74..84 shift down to 75..85.
$ws.Rows("74:74").Insert()

# Populate the newly inserted row 74 with the new weekly reading.
$ws.Range("A74").Value = 7
$ws.Range("B74").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C74").Value = "Ñuble"
$ws.Range("D74").Value = 44984
$ws.Range("D74").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E74").Value = 16
$ws.Range("F74").Value = "Fruta"
$ws.Range("G74").Value = 100103
$ws.Range("H74").Value = "Frutos de hueso (carozo)"
$ws.Range("I74").Value = 100103002
$ws.Range("J74").Value = "Ciruela"
$ws.Range("K74").Value = "Larry Ann"
$ws.Range("L74").Value = "Primera"
$ws.Range("M74").Value = 50
$ws.Range("N74").Value = 10000
$ws.Range("O74").Value = 10000
$ws.Range("P74").Value = 10000
$ws.Range("Q74").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R74").Value = "Región de O'Higgins"
$ws.Range("S74").Value = 556
$ws.Range("T74").Value = 18
